$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Update existing row 2 reference codes (A2:F2) to the new values
# ---------------------------------------------------------------
$ws.Range("A2").Value = "ME-638"
$ws.Range("B2").Value = "OF-643"
$ws.Range("C2").Value = "CT-308"
$ws.Range("D2").Value = "BR-434"
$ws.Range("E2").Value = "BI-146"
$ws.Range("F2").Value = "RE-76 "

# ---------------------------------------------------------------
# Add the two new process columns: G (Sale Order) and H (Purchase req)
# ---------------------------------------------------------------

# G1 header: copy header formatting from A1 (bold font, border, center
# alignment) then restyle the fill to the new grey header colour and
# set the caption text.
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Sale Order"
$ws.Range("G1").Interior.Color = 10921638

# G2 data cell: plain, unstyled (matches the other unstyled data cells).
$ws.Range("G2").Value = "OR-207"

# H1 header: same approach, with the gold/orange fill colour.
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Purchase req"
$ws.Range("H1").Interior.Color = 49407

# H2 data cell: copy the centered style used by F2.
$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = "POR-897"

# ---------------------------------------------------------------
# Column widths for the new columns (best-fit-like custom widths)
# ---------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 9.5
$ws.Columns("H").ColumnWidth = 11.65

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------
# Scroll / selection position
# ---------------------------------------------------------------
$ws.Range("H8").Select()
